# Rename Muni to Pilot
#
# The "Dictionary" sheet has an operator-name column (A) that lists
# "SF Muni" for a block of rows (795-1016). This renames every one of
# those cells to "SF Muni Pilot" (a whole-cell match/replace, so any
# unrelated strings such as "Free Muni for Seniors" are left alone),
# and then resets the sheet's view/selection to A2 to match the
# saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = $ws.Range("A1:A1016")
$colA.Replace("SF Muni", "SF Muni Pilot", 1, 1, $false, $false, $false)

$ws.Range("A2").Select()
